$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear old content that is being replaced / relocated ---
$ws.Range("A2:E4").ClearContents()

# --- Row 2: PhD in Neuroscience (Milena Vasquez-Amezquita) ---
$ws.Range("A2").Value = 'PhD in Neuroscience'
$ws.Range("B2").Value = '2015 - 2018'
$ws.Range("C2").Value = '\href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'
$ws.Range("D2").Value = '\href{https://www.uv.es/}{Universitat de València}, Spain'
$ws.Range("E2").Value = 'Thesis \textbf{(\textit{Summa Cum Laude})}: \textit{\href{http://hdl.handle.net/10550/67639}{Preferencias sexuales típicas y atípicas según sexo y edad de los estímulos: Utilidad de la técnica de rastreo ocular} [Typical and atypical sexual preferences according to sex and age of the stimuli: Usefulness of the eye tracking technique]}'

# --- Row 3: blank A-D, Alicia Salvador supervision note in E ---
$ws.Range("E3").Value = 'Supervised together with  Alicia Salvador'

# --- Row 4: Professional Doctorate (Francisco Javier Flores) ---
$ws.Range("A4").Value = 'Professional Doctorate in Counselling Psychology'
$ws.Range("B4").Value = '2015 - 2018'
$ws.Range("C4").Value = '\href{https://www.researchgate.net/profile/Francisco-Flores-14}{Francisco Javier Flores}'
$ws.Range("D4").Value = '\href{https://www.uel.ac.uk/}{U. of East London}, UK'
$ws.Range("E4").Value = 'Tésis: \textit{What sense do people make of the functions of their ’behaviours that may be causing problems in their everyday life’? A hybrid deductive/inductive template analysis}'

# --- Row 5: blank A-D, Lisa Chiara Fellin supervision note in E ---
$ws.Range("E5").Value = 'Supervised together with Lisa Chiara Fellin'

# --- Formatting: copy the existing wrap/left/top cell style onto the whole block ---
$ws.Range("C2").Copy() | Out-Null
$ws.Range("A2:E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(5).AutoFit()

# --- Selection ---
$ws.Range("A6:XFD7").Select() | Out-Null
